$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (pushes old rows 3-14 down to 4-15)
$ws.Rows("3:3").Insert()

# Fill new staff row (row 3)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "못생긴오리"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1000
$ws.Range("G3").Value = 3000
$ws.Range("H3").Value = 3001
$ws.Range("H3").VerticalAlignment = -4108

# Re-number Id column for the rest (A4:A15 = 2..13)
for ($i = 4; $i -le 15; $i++) {
    $ws.Range("A$i").Value = $i - 2
}

# Add new column I "PathFile" / "int32" header rows + data
$ws.Range("I1").Value = "PathFile"
$ws.Range("I2").Value = "int32"
for ($i = 3; $i -le 15; $i++) {
    $ws.Range("I$i").Value = 9001
}

# Match the recorded selection state (row 3 selected)
[void]$ws.Rows("3:3").Select()
